$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-27 Saturday" "2024-04-28 Sunday"

Replace-Text "20×85=1700" "40×29=1160"
Replace-Text "27×93=2511" "25×22=550"
Replace-Text "28×12=336" "74×59=4366"
Replace-Text "33×27=891" "72×23=1656"
Replace-Text "64×85=5440" "11×20=220"
Replace-Text "81×45=3645" "59×12=708"
Replace-Text "99×55=5445" "23×72=1656"
Replace-Text "55×19=1045" "18×60=1080"
Replace-Text "70×84=5880" "57×60=3420"
Replace-Text "12×82=984" "22×80=1760"
Replace-Text "87×36=3132" "53×73=3869"
Replace-Text "69×41=2829" "41×82=3362"
Replace-Text "44×12=528" "24×29=696"
Replace-Text "41×65=2665" "14×90=1260"
Replace-Text "52×86=4472" "60×83=4980"
Replace-Text "37×60=2220" "75×99=7425"
Replace-Text "45×42=1890" "24×96=2304"
Replace-Text "11×37=407" "67×61=4087"
Replace-Text "30×51=1530" "25×80=2000"
Replace-Text "95×72=6840" "69×53=3657"
Replace-Text "14×25=350" "15×13=195"
Replace-Text "85×69=5865" "94×24=2256"
Replace-Text "88×51=4488" "32×48=1536"
Replace-Text "82×59=4838" "51×60=3060"
Replace-Text "75×29=2175" "78×21=1638"

Write-Output "Done"
